$d = $word.ActiveDocument

$replacements = @(
    @{old="431÷5="; new="684÷2="},
    @{old="311÷8="; new="480÷4="},
    @{old="239÷3="; new="950÷8="},
    @{old="373÷2="; new="607÷3="},
    @{old="699÷3="; new="328÷5="},
    @{old="154÷7="; new="972÷2="},
    @{old="634÷2="; new="149÷9="},
    @{old="680÷7="; new="662÷8="},
    @{old="380÷2="; new="826÷8="},
    @{old="272÷8="; new="254÷3="},
    @{old="470÷4="; new="917÷4="},
    @{old="430÷2="; new="810÷5="},
    @{old="576÷7="; new="892÷8="},
    @{old="949÷9="; new="668÷9="},
    @{old="487÷7="; new="196÷4="},
    @{old="472÷6="; new="261÷2="},
    @{old="746÷2="; new="910÷9="},
    @{old="168÷4="; new="556÷7="},
    @{old="774÷8="; new="482÷6="},
    @{old="572÷6="; new="448÷3="},
    @{old="997÷9="; new="524÷2="},
    @{old="750÷4="; new="481÷9="},
    @{old="503÷4="; new="565÷4="},
    @{old="408÷8="; new="621÷5="},
    @{old="652÷8="; new="404÷9="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
